# Updated remaining queries for C3DC
# Replaces the old "id"-based join conditions with the new
# "study_id"/"participant_id"-based join conditions across every SQL
# query stored in column B (TabQuery) and C2 (StatQuery), and widens
# column C to fit the updated (longer) query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldJoin = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newJoin = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

$cellsToUpdate = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellsToUpdate) {
    $cell = $ws.Range($addr)
    $text = $cell.Value()
    if ($text -ne $null -and $text.Contains($oldJoin)) {
        $cell.Value = $text.Replace($oldJoin, $newJoin)
    }
}

# Widen column C now that the StatQuery text is longer, and drop the
# stale bestFit auto-sizing flag (Excel clears it on explicit resize).
$ws.Columns.Item(3).ColumnWidth = 66.6667
